# Update automatico via Actualizar 02-20-2021 12-50-03
# This mirrors a scheduled "availability checker" refresh: a brand-new
# timestamp is recorded for the most recent check (rows 2-15), while the
# previous batches of timestamps shift down one block of 14 rows
# (old rows 2-15 values -> rows 16-29, old rows 16-29 values -> rows 30-43).
# The oldest block (previously in rows 30-43) simply rolls off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newest = 44247.53459866466
$middle = 44247.51330931713
$oldest = 44247.49201875

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = $newest
}

for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = $middle
}

for ($r = 30; $r -le 43; $r++) {
    $ws.Cells.Item($r, 4).Value = $oldest
}
